$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (IP) on row 2: server's own IP changes from 192.168.1.113 to 127.0.0.1
$ws.Range("C2").Value = "127.0.0.1"

# Column E (SqlIP) on row 2: value stays 192.168.0.24, but gets the "Text" number
# format (same format already used by C2/B2/A2 - style index 1 in the original file)
$ws.Range("E2").NumberFormat = "@"

# Active cell / selection moves from C2 to E2
[void]$ws.Range("E2").Select()
